$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Bump the auto-date footer text from "2021. 03. 23." to "2021. 03. 24."
#    on every slide layout / the slide master / the handout master / the
#    notes master that carries a date placeholder.
# ---------------------------------------------------------------------------
$newDate = "2021. 03. 24."

$m = $p.SlideMaster
$m.Shapes.Item(4).TextFrame.TextRange.Text = $newDate

$dateShapeByLayout = @{
    2  = 3
    4  = 4
    5  = 6
    6  = 2
    7  = 1
    8  = 5
    9  = 5
    10 = 3
    11 = 3
}

foreach ($layoutIndex in $dateShapeByLayout.Keys) {
    $cl = $m.CustomLayouts.Item($layoutIndex)
    $shapeIndex = $dateShapeByLayout[$layoutIndex]
    $cl.Shapes.Item($shapeIndex).TextFrame.TextRange.Text = $newDate
}

$p.HandoutMaster.Shapes.Item(2).TextFrame.TextRange.Text = $newDate
$p.NotesMaster.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# ---------------------------------------------------------------------------
# 2. Remove the "Adient - INTERNAL" MSIPCMContentMarking footer shape that was
#    drawn directly on the slide master.
# ---------------------------------------------------------------------------
for ($i = $m.Shapes.Count; $i -ge 1; $i--) {
    if ($m.Shapes.Item($i).Name -eq "MSIPCMContentMarking") {
        $m.Shapes.Item($i).Delete()
    }
}

# ---------------------------------------------------------------------------
# 3. Reposition the picture on the "Technológiák" slide (slide 6).
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $sh = $s6.Shapes.Item($i)
    if ($sh.Name -eq "Kép 3") {
        $sh.Left = 5950396 / 12700
        $sh.Top = 3014945 / 12700
    }
}

# ---------------------------------------------------------------------------
# 4. "Hátralévő feladatok" slide (slide 8): re-title it, grow the title box,
#    and shrink/move the content placeholder to make room.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$title = $s8.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Főbb feladatok amin dolgozunk"
$title.Height = (1552600 / 12700) + 0.00001

$content = $s8.Shapes.Item(2)
$content.Top = 2420888 / 12700
$content.Height = 3751312 / 12700
